$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns -----------------------------------------------
# 1) a new column F "CopiaConsulentiDa" (forced-copy-of-consultants feature)
# 2) a new column H "LOB" (after the former F/"ProjectType" which becomes G)
$ws.Columns("F").Insert()
$ws.Columns("H").Insert()

# --- Column F (CopiaConsulentiDa): header then data ------------------------
# Only row 4 carries a value in the sample data (rows 2 and 3 stay blank)
$ws.Range("F1").Value = "CopiaConsulentiDa"
$ws.Range("F4").Value = "CIR_006_00"

# --- Column H (LOB): header then data ---------------------------------------
$ws.Range("H1").Value = "LOB"
$ws.Range("H2").Value = "ERP"
$ws.Range("H3").Value = "CDG"
$ws.Range("H4").Value = "SCP"

# --- Column widths -----------------------------------------------------------
$ws.Columns("F").ColumnWidth = 26.166666666666668   # -> stored width 27
$ws.Columns("H").ColumnWidth = 15                    # -> stored width ~15.83 (closest to former F width)

# --- Selection (bug fix for ore/spese input moved the active cell) ---------
[void]$ws.Range("H4").Select()
